# Maquilador, back end, bug de indicaciones
#
# Renames the generic "Indicaciones" template sheet/defined-name to
# "Maquilador" and reworks the listado columns from
# (Clave, Nombre, Descripcion, Activo) to
# (Clave, Nombre, Direccion, Correo, Telefono, Activo).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet / defined name rename -------------------------------------------
$ws.Name = "Maquilador"
$wb.Names.Item("Indicaciones").Delete()
$wb.Names.Add("Maquilador", "=Maquilador!`$A`$4:`$F`$5")

# --- Column contents ---------------------------------------------------
# Row 3 = headers, Row 4 = the mustache placeholders used by the report
# engine. A3 "Clave" / B3 "Nombre" and A4/B4 placeholders are untouched.
$ws.Range("F3").Value = "Activo"
$ws.Range("D3").Value = "Correo"

$ws.Range("A4").Value = "{{item.Clave}}"
$ws.Range("B4").Value = "{{item.Nombre}}"
$ws.Range("C4").Value = "{{item.Direccion}}"
$ws.Range("D4").Value = "{{item.Correo}}"
$ws.Range("E4").Value = "{{item.Telefono}}"
$ws.Range("F4").Value = "{{item.Activo}}"

$ws.Range("E3").Value = "Teléfono"
$ws.Range("C3").Value = "Dirección"

# Match the bold/centred header style (xf 1) on the new E3 header cell, and
# give the whole data row (row 4) the centred (non-bold) style used by the
# rest of the listado values.
$ws.Range("A3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("A4:F4").HorizontalAlignment = -4108

# --- Title banner merge now spans the two extra columns --------------------
$ws.Range("A1:D1").UnMerge()
$ws.Range("A1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$ws.Range("A1:F1").Merge()

# --- Column widths for the new Direccion / Correo / Telefono columns -------
$ws.Columns.Item(5).ColumnWidth = 24.28515625
$ws.Columns.Item(6).ColumnWidth = 24.140625

# --- Selection follows the new merged title range ---------------------------
$ws.Range("A1:F1").Select()
